# Textbox response formatting fix
# Renames task order sheets and updates stim/condition filenames with new timestamps.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (new timestamped identifiers) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-165116873726659"
$wb.Worksheets.Item(2).Name = "NB_TO-1651168739347002"
$wb.Worksheets.Item(3).Name = "RS_TO-16511687393480043"
$wb.Worksheets.Item(4).Name = "TOL_TO-1651168739393755"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687394695187"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511687372242396.csv"
$ws1.Range("B3").Value = "GNG_stims-1651168737250403.csv"
$ws1.Range("B4").Value = "go_stims-16511687372523959.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687372655983.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16511687383405037.csv"
$ws2.Range("B3").Value = "ZB-match_1-1651168737695932.csv"
$ws2.Range("B4").Value = "TB-16511687393317022.csv"
$ws2.Range("B5").Value = "OB-16511687383872411.csv"
$ws2.Range("B6").Value = "ZB-match_8-16511687373251407.csv"
$ws2.Range("B7").Value = "TB-16511687391160917.csv"
$ws2.Range("B8").Value = "ZB-match_0-16511687374728534.csv"
$ws2.Range("B9").Value = "TB-1651168738474611.csv"
$ws2.Range("B10").Value = "OB-1651168738030882.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687393621504.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687393490033.csv"
$ws4.Range("B4").Value = "MM_stims-16511687393769486.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687393621504.csv"
$ws4.Range("B6").Value = "MM_stims-16511687393927898.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687393769486.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16511687394534616.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511687394384012.csv"
$ws5.Range("B4").Value = "SAT_stims-16511687393978286.csv"
$ws5.Range("B5").Value = "SAT_stims-16511687394236147.csv"
